$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of an existing header cell (A1) so the new headers
# match the look (bold font, border, centered) of the other headers.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every player row.
for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = 89
    $ws.Cells.Item($row, 31).Value = 73
    $ws.Cells.Item($row, 32).Value = 0
}
